$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("B2").Value = 0.1374045801526718
$ws.Range("C2").Value = 0.6908396946564885
$ws.Range("J2").Value = 0.01526717557251908
$ws.Range("P2").Value = 0.09923664122137404
$ws.Range("S2").Value = 0.05725190839694656
$ws.Range("C3").Value = 0.02173913043478261
$ws.Range("J3").Value = 0.005434782608695652
$ws.Range("P3").Value = 0.8260869565217391
$ws.Range("S3").Value = 0.1467391304347826
$ws.Range("J4").Value = 0.0196078431372549
$ws.Range("P4").Value = 0.7843137254901961
$ws.Range("S4").Value = 0.196078431372549
$ws.Range("B6").Value = 0.07003891050583658
$ws.Range("D6").Value = 0.01556420233463035
$ws.Range("F6").Value = 0.07782101167315175
$ws.Range("J6").Value = 0.2490272373540856
$ws.Range("O6").Value = 0.01945525291828794
$ws.Range("Q6").Value = 0.1828793774319066
$ws.Range("R6").Value = 0.07003891050583658
$ws.Range("S6").Value = 0.3151750972762646
$ws.Range("B7").Value = 0.1170212765957447
$ws.Range("D7").Value = 0.03191489361702127
$ws.Range("F7").Value = 0.06914893617021277
$ws.Range("J7").Value = 0.0851063829787234
$ws.Range("O7").Value = 0.01595744680851064
$ws.Range("Q7").Value = 0.2553191489361702
$ws.Range("R7").Value = 0.09042553191489362
$ws.Range("S7").Value = 0.3351063829787234
$ws.Range("B8").Value = 0.09533898305084745
$ws.Range("D8").Value = 0.01271186440677966
$ws.Range("F8").Value = 0.06991525423728813
$ws.Range("J8").Value = 0.1059322033898305
$ws.Range("O8").Value = 0.02754237288135593
$ws.Range("Q8").Value = 0.1906779661016949
$ws.Range("R8").Value = 0.125
$ws.Range("S8").Value = 0.3728813559322034
$ws.Range("B9").Value = 0.06741573033707865
$ws.Range("D9").Value = 0.01685393258426966
$ws.Range("F9").Value = 0.1067415730337079
$ws.Range("J9").Value = 0.08426966292134831
$ws.Range("O9").Value = 0.01685393258426966
$ws.Range("Q9").Value = 0.1853932584269663
$ws.Range("R9").Value = 0.07865168539325842
$ws.Range("S9").Value = 0.4438202247191011
$ws.Range("B10").Value = 0.103363412633306
$ws.Range("D10").Value = 0.02625102543068089
$ws.Range("F10").Value = 0.07629204265791632
$ws.Range("J10").Value = 0.08941755537325677
$ws.Range("O10").Value = 0.01394585726004922
$ws.Range("Q10").Value = 0.2206726825266612
$ws.Range("R10").Value = 0.0992616899097621
$ws.Range("S10").Value = 0.3707957342083675
$ws.Range("G11").Value = 0.1580882352941176
$ws.Range("J11").Value = 0.08823529411764706
$ws.Range("K11").Value = 0.2058823529411765
$ws.Range("L11").Value = 0.5441176470588235
$ws.Range("S11").Value = 0.003676470588235294
$ws.Range("G12").Value = 0.7583892617449665
$ws.Range("J12").Value = 0.1879194630872483
$ws.Range("L12").Value = 0.01342281879194631
$ws.Range("S12").Value = 0.04026845637583892
$ws.Range("F15").Value = 0.03317535545023697
$ws.Range("H15").Value = 0.1042654028436019
$ws.Range("I15").Value = 0.04265402843601896
$ws.Range("J15").Value = 0.3933649289099526
$ws.Range("K15").Value = 0.07582938388625593
$ws.Range("M15").Value = 0.009478672985781991
$ws.Range("N15").Value = 0.004739336492890996
$ws.Range("O15").Value = 0.08056872037914692
$ws.Range("S15").Value = 0.2559241706161137
$ws.Range("F16").Value = 0.009708737864077669
$ws.Range("H16").Value = 0.2038834951456311
$ws.Range("I16").Value = 0.07281553398058252
$ws.Range("J16").Value = 0.4271844660194175
$ws.Range("K16").Value = 0.1019417475728155
$ws.Range("M16").Value = 0.01456310679611651
$ws.Range("O16").Value = 0.07281553398058252
$ws.Range("S16").Value = 0.0970873786407767
$ws.Range("F17").Value = 0.01649484536082474
$ws.Range("H17").Value = 0.1814432989690722
$ws.Range("I17").Value = 0.09484536082474226
$ws.Range("J17").Value = 0.4288659793814433
$ws.Range("K17").Value = 0.09484536082474226
$ws.Range("M17").Value = 0.01649484536082474
$ws.Range("N17").Value = 0.004123711340206186
$ws.Range("O17").Value = 0.0577319587628866
$ws.Range("S17").Value = 0.1051546391752577
$ws.Range("F18").Value = 0.03043478260869565
$ws.Range("H18").Value = 0.2217391304347826
$ws.Range("I18").Value = 0.08260869565217391
$ws.Range("J18").Value = 0.3956521739130435
$ws.Range("K18").Value = 0.1043478260869565
$ws.Range("M18").Value = 0.008695652173913044
$ws.Range("O18").Value = 0.05217391304347826
$ws.Range("S18").Value = 0.1043478260869565
$ws.Range("F19").Value = 0.02244389027431421
$ws.Range("H19").Value = 0.227763923524522
$ws.Range("I19").Value = 0.07564422277639235
$ws.Range("J19").Value = 0.3657522859517872
$ws.Range("K19").Value = 0.08894430590191189
$ws.Range("M19").Value = 0.02327514546965919
$ws.Range("N19").Value = 0.001662510390689942
$ws.Range("O19").Value = 0.0598503740648379
$ws.Range("S19").Value = 0.1346633416458853
